$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 92: 四方坪站
$ws.Cells.Item(92, 1).Value = 45946
$ws.Cells.Item(92, 2).Value = "四方坪站"
$ws.Cells.Item(92, 3).Formula = "=15159/126"
$ws.Cells.Item(92, 4).Formula = "=C92/(24*60)"
$ws.Cells.Item(92, 5).Formula = "=8637.61/126"
$ws.Cells.Item(92, 6).Formula = "=2985.93/126"
$ws.Cells.Item(92, 7).Formula = "=8637.61/(15159/60)"
$ws.Cells.Item(92, 8).Formula = "=374/126"

# Row 93: 高岭站
$ws.Cells.Item(93, 1).Value = 45946
$ws.Cells.Item(93, 2).Value = "高岭站"
$ws.Cells.Item(93, 3).Formula = "=5657/36"
$ws.Cells.Item(93, 4).Formula = "=C93/(24*60)"
$ws.Cells.Item(93, 5).Formula = "=3972.73/36"
$ws.Cells.Item(93, 6).Formula = "=1042.71/36"
$ws.Cells.Item(93, 7).Formula = "=3972.73/(5667/60)"
$ws.Cells.Item(93, 8).Formula = "=146/36"

$ws.Range("I94").Select() | Out-Null

Write-Host "done"
